# "update solution space size SAM to code in result files"
#
# The Projects sheet holds two tables: "Number of Trace Links in Gold
# Standard" (rows 4-8) and "Confusion matrix sums / Solution Space Size"
# (rows 13-17). Column E of the second table is "SAM-Code" — update its
# solution-space-size values for BigBlueButton, JabRef and TEAMMATES.

$wb = $excel.ActiveWorkbook

$sadCode  = $wb.Worksheets.Item("SAD-Code")
$projects = $wb.Worksheets.Item("Projects")

# --- data edits: SAM-Code (col E) solution space sizes on "Projects" ---
$projects.Activate()
$projects.Range("E13").Value = 13128
$projects.Range("E14").Value = 11874
$projects.Range("E16").Value = 13312

# --- selection bookkeeping left behind by the edit (moved from the old
#     C13:C17 block onto the E13:E17 block that was actually touched) ---
$sadCode.Activate()
$sadCode.Range("E13:E17").Select() | Out-Null

$projects.Activate()
$projects.Range("E13:E17").Select() | Out-Null
